# Remove the standalone "Actes" (italic) paragraph that immediately follows
# the "ACT" (book-code) heading paragraph. The "ACT" paragraph itself, and
# everything after the "Actes" paragraph, is left untouched.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "ACT`r") {
        $actesPara = $p.Next()
        $actesPara.Range.Delete()
        break
    }
}
